$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Min/Max parameter values
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 10

$ws.Range("B3").Value = 0.4
$ws.Range("C3").Value = 0.6

$ws.Range("B4").Value = 5
$ws.Range("C4").Value = 10

$ws.Range("B5").Value = 0.4
$ws.Range("C5").Value = 0.6

$ws.Range("B6").Value = 3
$ws.Range("C6").Value = 5

$ws.Range("B7").Value = 0.1
$ws.Range("C7").Value = 0.3

$ws.Range("B8").Value = 100
$ws.Range("C8").Value = 200

$ws.Range("C9").Value = 20

$ws.Range("C10").Value = 20

$ws.Range("C11").Value = 20

$ws.Range("B12").Value = 0.0001
$ws.Range("C12").Value = 20

# Update view: zoom and selection
[void]$ws.Range("C14").Select()
$excel.ActiveWindow.Zoom = 145
